# Add a new worksheet "ODI Batting Extra" right after "ODI Batting",
# matching the target workbook structure, and populate it with data.

$wb = $excel.ActiveWorkbook
$odiBatting = $wb.Worksheets.Item("ODI Batting")

$ws = $wb.Worksheets.Add($null, $odiBatting)
$ws.Name = "ODI Batting Extra"

# --- Header row (bold / bordered / centered, matching the style used by
#     the header rows on the other sheets) ---
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

$odiBatting.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# --- Data rows ---
# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$ws.Range("A2").Value = "'4452"
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = "'2"
$ws.Range("D2").Value = "'0"
$ws.Range("E2").Value = "'8.33%"
$ws.Range("F2").Value = "NO"

$ws.Range("A3").Value = "'4453"
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = "'0"
$ws.Range("D3").Value = "'0"
$ws.Range("E3").Value = "'0.36%"
$ws.Range("F3").Value = "NO"

$ws.Range("A4").Value = "'4563"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "'8"
$ws.Range("D4").Value = "'3"
$ws.Range("E4").Value = "'50.49%"
$ws.Range("F4").Value = "YES"

$ws.Range("A5").Value = "'4566"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "'0"
$ws.Range("D5").Value = "'0"
$ws.Range("E5").Value = "'0.38%"
$ws.Range("F5").Value = "NO"

$ws.Range("A6").Value = "'4568"
$ws.Range("F6").Value = "NO"

$ws.Range("A7").Value = "'4605"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "'0"
$ws.Range("D7").Value = "'0"
$ws.Range("E7").Value = "'0.33%"
$ws.Range("F7").Value = "NO"

$ws.Range("A8").Value = "'4608"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "'0"
$ws.Range("D8").Value = "'0"
$ws.Range("F8").Value = "NO"

$ws.Range("A9").Value = "'4614"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "'0"
$ws.Range("D9").Value = "'0"
$ws.Range("E9").Value = "'0.83%"
$ws.Range("F9").Value = "NO"

$ws.Range("A10").Value = "'4735"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "'2"
$ws.Range("D10").Value = "'1"
$ws.Range("E10").Value = "'9.49%"
$ws.Range("F10").Value = "NO"

$ws.Range("A11").Value = "'4745"
$ws.Range("F11").Value = "NO"
